$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 160.9723563333334
$ws.Range("H2").Value = 482.917069
$ws.Range("I2").Value = 0.3931645655589854
$ws.Range("J2").Value = 0.3931645655589854
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 46.71246076946579
$ws.Range("R2").Value = 420.412146925192
$ws.Range("S2").Value = 0.01348314463616053
$ws.Range("T2").Value = 0.01348314463616053
$ws.Range("G3").Value = 160.9723563333334
$ws.Range("H3").Value = 482.917069
$ws.Range("I3").Value = 0.3931645655589854
$ws.Range("J3").Value = 0.3931645655589854
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 1144.569901169621
$ws.Range("R3").Value = 10301.12911052659
$ws.Range("S3").Value = 0.3303701254324318
$ws.Range("T3").Value = 0.3303701254324318
$ws.Range("G4").Value = 160.9723563333334
$ws.Range("H4").Value = 482.917069
$ws.Range("I4").Value = 0.3931645655589854
$ws.Range("J4").Value = 0.3931645655589854
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 170.8393715445935
$ws.Range("R4").Value = 1537.554343901341
$ws.Range("S4").Value = 0.04931129549039307
$ws.Range("T4").Value = 0.04931129549039306
$ws.Range("I5").Value = 0.2197635343237224
$ws.Range("J5").Value = 0.2197635343237224
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 26.11042900333778
$ws.Range("R5").Value = 234.99386103004
$ws.Range("S5").Value = 0.007536547742617032
$ws.Range("T5").Value = 0.007536547742617033
$ws.Range("I6").Value = 0.2197635343237224
$ws.Range("J6").Value = 0.2197635343237224
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.1846639111456503
$ws.Range("T6").Value = 0.1846639111456504
$ws.Range("I7").Value = 0.2197635343237224
$ws.Range("J7").Value = 0.2197635343237224
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 95.49249190069945
$ws.Range("R7").Value = 859.432427106295
$ws.Range("S7").Value = 0.02756307543545502
$ws.Range("T7").Value = 0.02756307543545502
$ws.Range("I8").Value = 0.3870719001172923
$ws.Range("J8").Value = 0.3870719001172923
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 45.98858221997867
$ws.Range("R8").Value = 413.897239979808
$ws.Range("S8").Value = 0.01327420340247308
$ws.Range("T8").Value = 0.01327420340247308
$ws.Range("I9").Value = 0.3870719001172923
$ws.Range("J9").Value = 0.3870719001172923
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.3252505525550332
$ws.Range("T9").Value = 0.3252505525550332
$ws.Range("I10").Value = 0.3870719001172923
$ws.Range("J10").Value = 0.3870719001172923
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 168.1919632421427
$ws.Range("S10").Value = 0.04854714415978603
$ws.Range("T10").Value = 0.04854714415978603
